$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format temporarily so numeric-looking values
# (e.g. "309.72", "0.4764") are stored as text, matching the source data
# which uses inline strings throughout.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "27.038.54"
$ws.Cells.Item(2, 5).Value = "  +0.23%  "
$ws.Cells.Item(3, 4).Value = "1.847.95"
$ws.Cells.Item(3, 5).Value = "  +0.31%  "
$ws.Cells.Item(4, 5).Value = "  +0.83%  "
$ws.Cells.Item(5, 5).Value = "  +0.57%  "
$ws.Cells.Item(6, 4).Value = "309.72"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "
$ws.Cells.Item(7, 4).Value = "0.4764"
$ws.Cells.Item(8, 4).Value = "0.3685"
$ws.Cells.Item(8, 5).Value = "  +1.78%  "
$ws.Cells.Item(9, 4).Value = "0.07239"
$ws.Cells.Item(9, 5).Value = "  +1.53%  "
$ws.Cells.Item(10, 4).Value = "0.9316"
$ws.Cells.Item(10, 5).Value = "  +2.06%  "
$ws.Cells.Item(11, 5).Value = "  +1.71%  "
$ws.Cells.Item(12, 4).Value = "0.07785"
$ws.Cells.Item(12, 5).Value = "  +1.21%  "
$ws.Cells.Item(13, 4).Value = "1.872.42"
$ws.Cells.Item(13, 5).Value = "  +2.38%  "
$ws.Cells.Item(14, 5).Value = "  +2.30%  "
$ws.Cells.Item(15, 4).Value = "6.481"
$ws.Cells.Item(15, 5).Value = "  +1.24%  "
$ws.Cells.Item(16, 4).Value = "88.96"
$ws.Cells.Item(16, 5).Value = "  +0.79%  "
$ws.Cells.Item(17, 4).Value = "1.018"
$ws.Cells.Item(17, 5).Value = "  +0.77%  "
$ws.Cells.Item(18, 4).Value = "0.000008659"
$ws.Cells.Item(18, 5).Value = "  +0.88%  "
$ws.Cells.Item(19, 4).Value = "1.015"
$ws.Cells.Item(19, 5).Value = "  +0.69%  "
$ws.Cells.Item(20, 4).Value = "27.054.65"
$ws.Cells.Item(20, 5).Value = "  +0.11%  "
$ws.Cells.Item(21, 4).Value = "14.54"
$ws.Cells.Item(21, 5).Value = "  +1.48%  "
$ws.Cells.Item(22, 4).Value = "5.053"
$ws.Cells.Item(22, 5).Value = "  +0.73%  "
$ws.Cells.Item(23, 5).Value = "  +0.13%  "
$ws.Cells.Item(24, 4).Value = "1.930"
$ws.Cells.Item(24, 5).Value = "  +0.10%  "
$ws.Cells.Item(25, 4).Value = "152.80"
$ws.Cells.Item(25, 5).Value = "  +0.15%  "
$ws.Cells.Item(26, 4).Value = "18.36"
$ws.Cells.Item(26, 5).Value = "  +0.83%  "
$ws.Cells.Item(27, 4).Value = "1.990"
$ws.Cells.Item(27, 5).Value = "  -2.21%  "
$ws.Cells.Item(28, 4).Value = "114.64"
$ws.Cells.Item(28, 5).Value = "  +0.53%  "
$ws.Cells.Item(29, 4).Value = "4.932"
$ws.Cells.Item(29, 5).Value = "  +0.76%  "
$ws.Cells.Item(30, 4).Value = "0.08870"
$ws.Cells.Item(30, 5).Value = "  +0.16%  "
$ws.Cells.Item(31, 4).Value = "3.322"
$ws.Cells.Item(31, 5).Value = "  +3.88%  "
$ws.Cells.Item(32, 4).Value = "1.178"
$ws.Cells.Item(32, 5).Value = "  +0.56%  "
$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(33, 4).Value = "4.516"
$ws.Cells.Item(33, 5).Value = "  +1.26%  "
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(34, 4).Value = "0.7376"
$ws.Cells.Item(34, 5).Value = "  -1.40%  "
$ws.Cells.Item(35, 4).Value = "2.662"
$ws.Cells.Item(35, 5).Value = "  -6.28%  "
$ws.Cells.Item(36, 4).Value = "1.115"
$ws.Cells.Item(36, 5).Value = "  +3.17%  "
$ws.Cells.Item(37, 4).Value = "0.01971"
$ws.Cells.Item(37, 5).Value = "  +1.63%  "
$ws.Cells.Item(38, 4).Value = "0.05260"
$ws.Cells.Item(38, 5).Value = "  +1.84%  "
$ws.Cells.Item(39, 4).Value = "2.969"
$ws.Cells.Item(39, 5).Value = "  -0.56%  "
$ws.Cells.Item(40, 4).Value = "0.5284"
$ws.Cells.Item(40, 5).Value = "  +2.02%  "
$ws.Cells.Item(41, 4).Value = "7.029"
$ws.Cells.Item(41, 5).Value = "  +1.78%  "
$ws.Cells.Item(42, 4).Value = "0.1521"
$ws.Cells.Item(42, 5).Value = "  +0.56%  "
$ws.Cells.Item(43, 4).Value = "8.283"
$ws.Cells.Item(43, 5).Value = "  +1.96%  "
$ws.Cells.Item(44, 4).Value = "10.57"
$ws.Cells.Item(44, 5).Value = "  +0.54%  "
$ws.Cells.Item(45, 4).Value = "0.4744"
$ws.Cells.Item(45, 5).Value = "  +1.04%  "
$ws.Cells.Item(46, 4).Value = "1.016"
$ws.Cells.Item(46, 5).Value = "  +0.64%  "
$ws.Cells.Item(47, 4).Value = "101.87"
$ws.Cells.Item(47, 5).Value = "  +1.20%  "
$ws.Cells.Item(48, 4).Value = "1.614"
$ws.Cells.Item(48, 5).Value = "  +0.60%  "
$ws.Cells.Item(49, 4).Value = "65.87"
$ws.Cells.Item(49, 5).Value = "  +2.31%  "
$ws.Cells.Item(50, 4).Value = "0.06064"
$ws.Cells.Item(50, 5).Value = "  +0.34%  "
$ws.Cells.Item(51, 4).Value = "0.8934"
$ws.Cells.Item(51, 5).Value = "  +3.76%  "

# Restore column D to the default (Normal) style so no stray cell-level
# number formatting is left behind.
$dRange.Style = "Normal"
